$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44469
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 1200

$ws.Range("D3").Value = 44463
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("P3").Value = 1200

$ws.Range("D4").Value = 44698
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 1100

$ws.Range("D5").Value = 44656
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 1000

$ws.Range("D6").Value = 44369
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 8000
$ws.Range("P6").Value = 800

$ws.Range("D7").Value = 44525
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("P7").Value = 900

$ws.Range("D8").Value = 44707
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("P8").Value = 1200

$ws.Range("D9").Value = 44705
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 1200

$ws.Range("D10").Value = 44473
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("P10").Value = 1100

$ws.Range("D11").Value = 44645
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 1000

$ws.Range("D12").Value = 44530
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("P12").Value = 1000

$ws.Range("D13").Value = 44425
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 13000
$ws.Range("P13").Value = 1300

$ws.Range("D14").Value = 44348
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("P14").Value = 1000

$ws.Range("D15").Value = 44523
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("P15").Value = 900

$ws.Range("D16").Value = 44649
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("P16").Value = 1000

$ws.Range("D17").Value = 44663
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("P17").Value = 1200

$ws.Range("D18").Value = 44708
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 11000
$ws.Range("L18").Value = 11000
$ws.Range("M18").Value = 11000
$ws.Range("P18").Value = 1100

$ws.Range("D19").Value = 44659
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("P19").Value = 1000

$ws.Range("D20").Value = 44526
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 9000
$ws.Range("P20").Value = 900

